$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new column H ("p_adj") with header + values ------------------
$ws.Cells.Item(1, 8).Value = "p_adj"

$ws.Cells.Item(2, 8).Value  = 0.73879678133642601
$ws.Cells.Item(3, 8).Value  = 0.92080116715127303
$ws.Cells.Item(4, 8).Value  = 0.026509682071014701
$ws.Cells.Item(5, 8).Value  = 0.121016090533018
$ws.Cells.Item(6, 8).Value  = 0.024948821318035601
$ws.Cells.Item(7, 8).Value  = 0.032734395379120103
$ws.Cells.Item(8, 8).Value  = 0.73879678133642601
$ws.Cells.Item(9, 8).Value  = 1
$ws.Cells.Item(10, 8).Value = 0.026509682071014701
$ws.Cells.Item(11, 8).Value = 1
$ws.Cells.Item(12, 8).Value = 0.091555621782598004
$ws.Cells.Item(13, 8).Value = 0.026509682071014701
$ws.Cells.Item(14, 8).Value = 0.90906497772052197
$ws.Cells.Item(15, 8).Value = 0.11466435817746801
$ws.Cells.Item(16, 8).Value = 0.79844704799307997
$ws.Cells.Item(17, 8).Value = 0.73879678133642601
$ws.Cells.Item(18, 8).Value = 0.45194446046770398
$ws.Cells.Item(19, 8).Value = 0.046249852908051302

# --- Move the "significant p-value" yellow highlight from column E ----
# --- (raw p) to column H (adjusted p) for the rows that used to be ----
# --- flagged -----------------------------------------------------------
$sigRows = @(4, 6, 7, 10, 13, 19)
foreach ($r in $sigRows) {
    $ws.Cells.Item($r, 5).ClearFormats()
    $ws.Cells.Item($r, 8).Interior.Color = 65535
}

# --- Column B ("domain") was widened / best-fit by the author ---------
$ws.Columns.Item(2).ColumnWidth = 23.416667

# --- Update selection / scroll position --------------------------------
$ws.Range("H13").Select() | Out-Null
